$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: update title and link
$ws.Range("D12").Value = "“혼자 공부하는 머신러닝+딥러닝(개정판)” 코드 업데이트 안내"
$ws.Range("E12").Value = "https://tensorflow.blog/2025/06/09/%ed%98%bc%ec%9e%90-%ea%b3%b5%eb%b6%80%ed%95%98%eb%8a%94-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d%eb%94%a5%eb%9f%ac%eb%8b%9d%ea%b0%9c%ec%a0%95%ed%8c%90-%ec%bd%94%eb%93%9c-%ec%97%85%eb%8d%b0%ec%9d%b4/"

# Row 18: update title and link
$ws.Range("D18").Value = "여행의 겹쳐 읽기"
$ws.Range("E18").Value = "https://freesearch.pe.kr/archives/5533"

# Row 42: update title and link
$ws.Range("D42").Value = "[Window/Network]IP설정 - 윈도우 멀티 IP 설정"
$ws.Range("E42").Value = "https://kjk92.tistory.com/131"
